$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 2797
$ws.Range("F5").Value = 982
$ws.Range("F10").Value = 779
$ws.Range("F13").Value = 594
$ws.Range("F14").Value = 1178
$ws.Range("F16").Value = 719
$ws.Range("F17").Value = 640
$ws.Range("F22").Value = 776
$ws.Range("F23").Value = 8273
$ws.Range("F24").Value = 545
$ws.Range("F25").Value = 545
$ws.Range("F32").Value = 16
$ws.Range("F33").Value = 248
$ws.Range("F37").Value = 223
$ws.Range("F39").Value = 43
$ws.Range("F40").Value = 86
$ws.Range("F43").Value = 155

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 221
$ws.Range("F12").Value = 62
$ws.Range("F15").Value = 49
$ws.Range("F16").Value = 244

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 784

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 784
$ws.Range("F4").Value = 2797
$ws.Range("F5").Value = 982
$ws.Range("F10").Value = 779
$ws.Range("F15").Value = 594
$ws.Range("F16").Value = 1178
$ws.Range("F20").Value = 719
$ws.Range("F21").Value = 640
$ws.Range("F25").Value = 776
$ws.Range("F26").Value = 8273
$ws.Range("F27").Value = 221
$ws.Range("F28").Value = 545
$ws.Range("F29").Value = 545
$ws.Range("F34").Value = 16
$ws.Range("F35").Value = 248
$ws.Range("F38").Value = 62
$ws.Range("F39").Value = 62
$ws.Range("F41").Value = 49
$ws.Range("F42").Value = 223
$ws.Range("F45").Value = 43
$ws.Range("F46").Value = 86
